# Update the equilibrated-data values in B2:C16 with the recalculated
# AGW/ABW figures referenced by the commit ("Add ABW and AGW experimental
# data...").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(40.454047619047607, 0.51896763945578228)
    3  = @(38.326428571428558, 0.50536133673469386)
    4  = @(19.988571428571429, 0.36464865306122451)
    5  = @(37.770238095238092, 0.5023441666666667)
    6  = @(38.499523809523808, 0.50709372789115648)
    7  = @(28.39595238095238, 0.43526938435374152)
    8  = @(28.653571428571428, 0.43676229591836752)
    9  = @(28.26428571428572, 0.43405867346938781)
    10 = @(39.940714285714293, 0.51580579591836739)
    11 = @(31.722380952380949, 0.45952134693877539)
    12 = @(32.367142857142859, 0.46562446938775509)
    13 = @(32.304285714285712, 0.4637972448979592)
    14 = @(35.265000000000001, 0.48464185714285718)
    15 = @(35.556904761904768, 0.48712959523809529)
    16 = @(35.808809523809529, 0.48853447278911571)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
}

# Match the saved view state: last selected cell was N16.
$ws.Range("N16").Select()
